$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.503.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.23%  "
$ws.Range("D3").Value = "'1.841.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.77%  "
$ws.Range("D4").Value = "'1.029"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +2.79%  "
$ws.Range("D5").Value = "'319.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.25%  "
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("D7").Value = "'0.4370"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("D8").Value = "'0.3730"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").Value = "'0.07390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("D10").Value = "'0.8769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +4.82%  "
$ws.Range("D12").Value = "'1.857.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.37%  "
$ws.Range("D13").Value = "'5.494"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.67%  "
$ws.Range("D14").Value = "'6.671"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.61%  "
$ws.Range("D15").Value = "'0.07148"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").Value = "'82.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.67%  "
$ws.Range("D17").Value = "'1.032"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.99%  "
$ws.Range("D18").Value = "'0.000009027"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.25%  "
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("D21").Value = "'27.519.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.24%  "
$ws.Range("D22").Value = "'5.231"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("D24").Value = "'2.077.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.87%  "
$ws.Range("D25").Value = "'156.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").Value = "'1.922"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.58%  "
$ws.Range("D27").Value = "'18.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("D28").Value = "'5.254"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.27%  "
$ws.Range("D29").Value = "'1.936"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.36%  "
$ws.Range("D30").Value = "'116.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "'0.09063"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").Value = "'1.210"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.87%  "
$ws.Range("D33").Value = "'0.7626"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("D34").Value = "'4.485"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("D35").Value = "'2.876"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.24%  "
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'1.147"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("D38").Value = "'0.01971"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.45%  "
$ws.Range("D39").Value = "'0.05251"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("D40").Value = "'0.5173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.21%  "
$ws.Range("D41").Value = "'2.781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.02%  "
$ws.Range("E42").Value = "  +3.41%  "
$ws.Range("D43").Value = "'6.628"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.73%  "
$ws.Range("D44").Value = "'8.511"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.92%  "
$ws.Range("D45").Value = "'108.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("D46").Value = "'10.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "'1.707"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.95%  "
$ws.Range("D49").Value = "'0.4639"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.37%  "
$ws.Range("D50").Value = "'0.06331"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").Value = "'1.890"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.91%  "
